$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.472.33"
$ws.Range("E2").Value = "  +1.85%  "
$ws.Range("D3").Value = "1.688.61"
$ws.Range("E3").Value = "  +3.68%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "221.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.523"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "30.45"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.45%  "
$ws.Range("E9").Value = "  +2.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0624"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0902"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("D12").Value = "1.929.66"
$ws.Range("E12").Value = "  +3.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +15.90%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.701.83"
$ws.Range("E14").Value = "  +4.59%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.617"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.61%  "
$ws.Range("D17").Value = "30.527.84"
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").Value = "0.0₃0720"
$ws.Range("E20").Value = "  +2.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.112"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("E30").Value = "  +2.66%  "
$ws.Range("E31").Value = "  +4.44%  "
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("E33").Value = "  +3.65%  "
$ws.Range("D34").Value = "1.509.06"
$ws.Range("E34").Value = "  +5.58%  "
$ws.Range("E35").Value = "  +5.50%  "
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0179"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.93%  "
$ws.Range("E38").Value = "  -4.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "79.26"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.42%  "
$ws.Range("E40").Value = "  +5.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.22%  "
$ws.Range("E42").Value = "  +2.94%  "
$ws.Range("E43").Value = "  +1.95%  "
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.995"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "52.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.99%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.89%  "
$ws.Range("D48").Value = "1.821.03"
$ws.Range("E48").Value = "  +2.89%  "
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "95.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.42%  "
$ws.Range("D51").Value = "0.0₆0112"
$ws.Range("E51").Value = "  +4.72%  "
